# perbaikan daftar mahasiswa alpha
# Rebuild Sheet1 so it holds the corrected "alpha" attendance list:
# nim | alpha | poin | status | periode  (was: mahasiswa_id | sakit | izin | alpha | poin | status | periode + kompen cols)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- headers (row 1) ---
$ws.Range("A1").Value = "nim"
$ws.Range("B1").Value = "alpha"
$ws.Range("C1").Value = "poin"
$ws.Range("D1").Value = "status"
$ws.Range("E1").Value = "periode"

# --- data row 2 ---
$ws.Range("A2").Value = 2241760112
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "Lunas"
$ws.Range("E2").Value = "2024/2025"

# --- data row 3 ---
$ws.Range("A3").Value = 2241760111
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "Belum Lunas"
$ws.Range("E3").Value = "2024/2025"

# the old layout used columns F and G too (kompen / periode) - drop them entirely
$ws.Range("F1:G3").Clear()

# column widths tweaked by the author after the rework
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.083333333333334

# selection left on D4 after the edit
$ws.Range("D4").Select() | Out-Null
